$d = $word.ActiveDocument

# --- Edit 1: merge the three runs in the "Stern über Bethlehem" bullet
#     into a single run "Stern über Bethlehem – Sarah" ---
$p1 = $d.Paragraphs(2)
$r1 = $p1.Range
$r1.End = $r1.End - 1          # exclude the paragraph mark
$r1.Delete()
$r1.InsertAfter("Stern über Bethlehem – Sarah")

# --- Edit 2: remove the "?? – Jan-Hendrik" bullet that follows
#     "Alle Jahre wieder – Sarah" ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "?? – Jan-Hendrik" -and $i -gt 1) {
        $prevText = $d.Paragraphs($i - 1).Range.Text.TrimEnd([char]13, [char]7)
        if ($prevText -eq "Alle Jahre wieder – Sarah") {
            $targetIndex = $i
        }
    }
}

if ($targetIndex -gt 0) {
    $d.Paragraphs($targetIndex).Range.Delete()
}
